$d = $word.ActiveDocument

$d.Content.Find.Execute("35×40=", $true, $false, $false, $false, $false, $true, 1, $false, "34×13=", 2) | Out-Null
$d.Content.Find.Execute("43×99=", $true, $false, $false, $false, $false, $true, 1, $false, "57×15=", 2) | Out-Null
$d.Content.Find.Execute("49×95=", $true, $false, $false, $false, $false, $true, 1, $false, "51×30=", 2) | Out-Null
$d.Content.Find.Execute("82×74=", $true, $false, $false, $false, $false, $true, 1, $false, "17×62=", 2) | Out-Null
$d.Content.Find.Execute("94×26=", $true, $false, $false, $false, $false, $true, 1, $false, "29×18=", 2) | Out-Null
$d.Content.Find.Execute("54×35=", $true, $false, $false, $false, $false, $true, 1, $false, "36×32=", 2) | Out-Null
$d.Content.Find.Execute("43×76=", $true, $false, $false, $false, $false, $true, 1, $false, "21×69=", 2) | Out-Null
$d.Content.Find.Execute("87×57=", $true, $false, $false, $false, $false, $true, 1, $false, "55×70=", 2) | Out-Null
$d.Content.Find.Execute("13×91=", $true, $false, $false, $false, $false, $true, 1, $false, "51×42=", 2) | Out-Null
$d.Content.Find.Execute("92×90=", $true, $false, $false, $false, $false, $true, 1, $false, "27×25=", 2) | Out-Null
$d.Content.Find.Execute("99×89=", $true, $false, $false, $false, $false, $true, 1, $false, "49×84=", 2) | Out-Null
$d.Content.Find.Execute("80×42=", $true, $false, $false, $false, $false, $true, 1, $false, "18×95=", 2) | Out-Null
$d.Content.Find.Execute("38×41=", $true, $false, $false, $false, $false, $true, 1, $false, "21×79=", 2) | Out-Null
$d.Content.Find.Execute("15×26=", $true, $false, $false, $false, $false, $true, 1, $false, "74×92=", 2) | Out-Null
$d.Content.Find.Execute("26×86=", $true, $false, $false, $false, $false, $true, 1, $false, "38×61=", 2) | Out-Null
$d.Content.Find.Execute("44×26=", $true, $false, $false, $false, $false, $true, 1, $false, "60×53=", 2) | Out-Null
$d.Content.Find.Execute("55×18=", $true, $false, $false, $false, $false, $true, 1, $false, "77×23=", 2) | Out-Null
$d.Content.Find.Execute("94×49=", $true, $false, $false, $false, $false, $true, 1, $false, "24×47=", 2) | Out-Null
$d.Content.Find.Execute("61×77=", $true, $false, $false, $false, $false, $true, 1, $false, "24×37=", 2) | Out-Null
$d.Content.Find.Execute("60×94=", $true, $false, $false, $false, $false, $true, 1, $false, "56×82=", 2) | Out-Null
$d.Content.Find.Execute("67×94=", $true, $false, $false, $false, $false, $true, 1, $false, "37×11=", 2) | Out-Null
$d.Content.Find.Execute("80×41=", $true, $false, $false, $false, $false, $true, 1, $false, "81×77=", 2) | Out-Null
$d.Content.Find.Execute("68×73=", $true, $false, $false, $false, $false, $true, 1, $false, "24×54=", 2) | Out-Null
$d.Content.Find.Execute("31×77=", $true, $false, $false, $false, $false, $true, 1, $false, "66×66=", 2) | Out-Null
$d.Content.Find.Execute("96×52=", $true, $false, $false, $false, $false, $true, 1, $false, "77×93=", 2) | Out-Null
